# SR [2022-09-06] Refactoring in Process --> adding to PS Class Build, Launch, WatchMode
$wb = $excel.ActiveWorkbook

$wsPlanned = $wb.Worksheets.Item("Planned Objects")
$wsTeam    = $wb.Worksheets.Item("Team")

# --- "Team" sheet: clear the sample e-mail placeholder, update cursor position ---
$wsTeam.Activate()
$wsTeam.Range("D4").ClearContents()
$wsTeam.Range("E7").Select()

# --- "Planned Objects" sheet: update cursor position, keep it the active tab ---
$wsPlanned.Activate()
$wsPlanned.Range("D11").Select()

# --- Append a new blank worksheet named "Sheet1" at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Sheet1"

# Restore "Planned Objects" as the active sheet/tab after the insert.
$wsPlanned.Activate()
$wsPlanned.Range("D11").Select()
